# "Generate Report for Handback"
#
# The handback transform failed for the bb4a462f... file in both the
# zh-cn and de-de locales: the handback archive's inner file name
# ("punynbtu.4sw") didn't match the handoff file name it was supposed to
# correspond to. Update the status to reflect the failure and record the
# error detail message for each locale sheet, widening the "Error Detail"
# column so the message is readable.

$wb = $excel.ActiveWorkbook

# -- Overview sheet: the "zh-cn"/"de-de" status cells for the
#    bb4a462f... row both shared the "Ready for handoff" text; they now
#    read "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# -- zh-cn sheet: Status cell + new Error Detail text, and widen the
#    Error Detail column so the message fits.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsZhCn.Range("P3").Value = "Handback file name: punynbtu.4sw is different with handoff file name: bb4a462f-daab-4821-bfb7-7a196ea5815a.855705cc0274f0ea3bd67f4bde4796b15ae3f83a.zh-cn."

# -- de-de sheet: same treatment.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Range("P3").Value = "Handback file name: punynbtu.4sw is different with handoff file name: bb4a462f-daab-4821-bfb7-7a196ea5815a.855705cc0274f0ea3bd67f4bde4796b15ae3f83a.de-de."
